$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "B5M3N2BW0MA3"
$ws.Range("A3").Value = "XA43JHAGQ8V3"
$ws.Range("A4").Value = "TK50GFCXFHCN"
$ws.Range("A15").Value = "9M1A883VTX21"
$ws.Range("A16").Value = "D2ANGGG71FGC"
$ws.Range("A17").Value = "8AW7QA18SBTA"
$ws.Range("A18").Value = "EA5XZ049QR7S"
$ws.Range("A19").Value = "PTV3TSFPBF6W"
$ws.Range("A20").Value = "NHYK5008HQDA"
$ws.Range("A21").Value = "ZADGNDVPP03M"

$ws.Range("A12").Value = $null

$ws.Range("A7").Select()
